# Updates the cryptos list (Coin, Link, Price, Volume(1h)) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "67.107.98", "  -1.15%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.315.33", "  +1.88%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.01%  "),
    @(5, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "184.67", "  +0.34%  "),
    @(6, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "576.50", "  -0.81%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.00%  "),
    @(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.602", "  +0.57%  "),
    @(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.128", "  -0.38%  "),
    @(10, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "6.65", "  +0.79%  "),
    @(11, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.406", "  -0.07%  "),
    @(12, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.892.82", "  +1.85%  "),
    @(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.138", "  -0.70%  "),
    @(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "27.26", "  +0.09%  "),
    @(15, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "67.289.73", "  -0.92%  "),
    @(16, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000166", "  -0.40%  "),
    @(17, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.321.75", "  +2.05%  "),
    @(18, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "442.64", "  +6.94%  "),
    @(19, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.65", "  -0.67%  "),
    @(20, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "13.52", "  +2.37%  "),
    @(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.69", "  +2.47%  "),
    @(22, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "74.00", "  +4.06%  "),
    @(23, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.09%  "),
    @(24, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.464.16", "  +1.94%  "),
    @(25, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.511", "  +1.08%  "),
    @(26, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0000118", "  +1.26%  "),
    @(27, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.190", "  +1.70%  "),
    @(28, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "9.01", "  -3.75%  "),
    @(29, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  -0.19%  "),
    @(30, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.96", "  +1.17%  "),
    @(31, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "22.84", "  +1.40%  "),
    @(32, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "5.31", "  -2.08%  "),
    @(33, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "0.998", "  -0.05%  "),
    @(34, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.79", "  -0.51%  "),
    @(35, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.23", "  -0.43%  "),
    @(36, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "162.35", "  -0.39%  "),
    @(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.50", "  +4.53%  "),
    @(38, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.83", "  -2.04%  "),
    @(39, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "27.04", "  +0.81%  "),
    @(40, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.784.25", "  +5.77%  "),
    @(41, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.789", "  -0.35%  "),
    @(42, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.46", "  +0.67%  "),
    @(43, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "6.22", "  -1.12%  "),
    @(44, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "40.32", "  -1.03%  "),
    @(45, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0669", "  -0.34%  "),
    @(46, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "24.54", "  +1.62%  "),
    @(47, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.36", "  -1.74%  "),
    @(48, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "324.43", "  -3.72%  "),
    @(49, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0271", "  -0.53%  "),
    @(50, "ONDO", "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo", "0.982", "  +1.01%  "),
    @(51, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "6.15", "  -1.01%  ")
)

# Price (D) and Volume (E) columns contain numeric-looking text (e.g. "67.107.98",
# "1.00", "  -1.15%  ") that must stay plain text, exactly like the original
# inlineStr cells. Flip the range to Text format before writing so Excel does
# not auto-convert them to numbers, then restore the original cell style so the
# saved workbook doesn't pick up a stray number format.
$valueRange = $ws.Range("D2:E51")
$origStyle = $valueRange.Style

$valueRange.NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$valueRange.Style = $origStyle
